# Apply the edits described in the commit:
# "found sources for district heat and electricity use in PJ"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 18: "production|district heat|District heat use" was re-labelled
#    "consumption|district heat|District heat use" and real source values
#    (in PJ) were filled in for 2035 / 2050, replacing the placeholder
#    negative numbers.
# ---------------------------------------------------------------------
$ws.Range("E18").Value2 = "consumption|district heat|District heat use"
$ws.Range("G18").Value2 = 74.2
$ws.Range("H18").Value2 = 68.8

# Clear the red "needs checking" highlight fill that used to mark this row
# as suspect, now that a real source has been found (-4142 = xlColorIndexNone).
$ws.Range("E18").Interior.ColorIndex = -4142

# ---------------------------------------------------------------------
# 2) Row 21: add a newly-sourced variable - residential electricity use.
# ---------------------------------------------------------------------
$ws.Range("A21").Value2 = "remind"
$ws.Range("B21").Value2 = "SSP2-Base"
$ws.Range("C21").Value2 = "Wind"
$ws.Range("D21").Value2 = "DK"
$ws.Range("E21").Value2 = "consumption|electricity|residential electricity use"
$ws.Range("F21").Value2 = "Gwh"
$ws.Range("G21").Value2 = 8394
$ws.Range("H21").Value2 = 7724.71

# ---------------------------------------------------------------------
# 3) Selection moved to E18 (the cell that was just sourced/edited).
# ---------------------------------------------------------------------
$ws.Range("E18").Select()
